$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text format so that numeric-looking
# strings (e.g. "1.001", "93.00") are preserved verbatim as text, matching
# the inline-string cell type used throughout this sheet.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.226.47'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.820.37'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '313.19'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '0.4468'
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').Value = '0.3768'
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').Value = '0.07394'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = '0.8782'
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').Value = '20.81'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').Value = '1.822.13'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '6.707'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '5.413'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '93.00'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '0.07136'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '0.000008785'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '15.02'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '27.226.61'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').Value = '5.357'
$ws.Range('E22').Value = '  +3.79%  '
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '1.967'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = '151.07'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '2.290'
$ws.Range('E26').Value = '  +3.21%  '
$ws.Range('D27').Value = '18.55'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('D28').Value = '5.332'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('D29').Value = '117.29'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('D30').Value = '0.08863'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').Value = '0.7803'
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('D32').Value = '1.190'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').Value = '4.560'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('D35').Value = '1.002'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '1.111'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').Value = '0.01975'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').Value = '0.05253'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D39').Value = '7.321'
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('D40').Value = '0.5287'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '2.867'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '0.1703'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = '2.268'
$ws.Range('E43').Value = '  +14.81%  '
$ws.Range('D44').Value = '8.579'
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('D45').Value = '0.5026'
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('D46').Value = '10.54'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '104.75'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.684'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '1.001'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('D51').Value = '65.84'
$ws.Range('E51').Value = '  +4.63%  '

# Restore the default (unstyled) cell style now that the values are
# committed as text, so no stray number-format style lingers on the cells.
$dataRange.Style = "Normal"
